$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Return the raw OOXML for a single paragraph (the <w:p>...</w:p> fragment),
# extracted from Range.WordOpenXML so it is byte-for-byte faithful to the
# live document (no hand-transcription risk).
function Get-ParaXml($para) {
    $full = $para.Range.WordOpenXML
    $startTag = "<w:body>"
    $idx1 = $full.IndexOf($startTag) + $startTag.Length
    $idx2 = $full.IndexOf("</w:p>", $idx1) + "</w:p>".Length
    return $full.Substring($idx1, $idx2 - $idx1)
}

# Locate a paragraph by its exact (trimmed) text; returns the 1-based index
# or -1 if not found.
function Find-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Strip the attributes off the opening <w:p ...> tag, turning it into a bare
# <w:p> -- this is what freshly-inserted paragraphs look like.
function Make-BareP($xml) {
    $closeIdx = $xml.IndexOf(">")
    return "<w:p>" + $xml.Substring($closeIdx + 1)
}

# Replace the single <w:t>$oldText</w:t> run text with new text.
function Set-ParaText($xml, $oldText, $newText) {
    return $xml.Replace("<w:t>$oldText</w:t>", "<w:t>$newText</w:t>")
}

# Drop the run-level shading (<w:shd .../>) some of these bullets carry.
function Remove-Shd($xml) {
    return $xml.Replace('<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>', '')
}

# Insert a <w:lastRenderedPageBreak/> immediately before the given run text.
function Add-PageBreakBeforeText($xml, $text) {
    return $xml.Replace("<w:t>$text</w:t>", "<w:lastRenderedPageBreak/><w:t>$text</w:t>")
}

# Remove a <w:lastRenderedPageBreak/> element entirely.
function Remove-PageBreak($xml) {
    return $xml.Replace('<w:lastRenderedPageBreak/>', '')
}

# ---------------------------------------------------------------------------
# Locate the anchor paragraphs before any edits happen (so indices are
# still the original ones).
# ---------------------------------------------------------------------------

$idxSignup = Find-ParaIndexByText $d "Signup"
$idxForgotPw = Find-ParaIndexByText $d "Forgot Password"
$idxProfile = Find-ParaIndexByText $d "Profile"
$idxLogout = Find-ParaIndexByText $d "Logout"
$idxNonfunc = Find-ParaIndexByText $d "Nonfunctional Requirements"

$signupXml = Get-ParaXml($d.Paragraphs($idxSignup))
$forgotPwXml = Get-ParaXml($d.Paragraphs($idxForgotPw))
$profileXml = Get-ParaXml($d.Paragraphs($idxProfile))
$logoutXml = Get-ParaXml($d.Paragraphs($idxLogout))
$acctBalXml = Get-ParaXml($d.Paragraphs($idxLogout - 1))
$nonfuncXml = Get-ParaXml($d.Paragraphs($idxNonfunc))

# ---------------------------------------------------------------------------
# Build the replacement paragraphs.
# ---------------------------------------------------------------------------

# "Signup" -> "Register", losing its run shading.
$registerXml = Remove-Shd (Set-ParaText $signupXml "Signup" "Register")

# New paragraph "Forgot Password" (clone of the bullet template, keeps shading).
$newForgotPwXml = Set-ParaText (Make-BareP $signupXml) "Signup" "Forgot Password"

# Old "Forgot Password" paragraph -> "Profile" (keeps its own identity/shading).
$newProfileXml = Set-ParaText $forgotPwXml "Forgot Password" "Profile"

# Old "Profile" paragraph -> "Transactions".
$newTransactionsXml = Set-ParaText $profileXml "Profile" "Transactions"

# New paragraphs "Deposit" and "Withdrawal" (clones of the bullet template).
$newDepositXml = Set-ParaText (Make-BareP $signupXml) "Signup" "Deposit"
$newWithdrawalXml = Set-ParaText (Make-BareP $signupXml) "Signup" "Withdrawal"

# "Account Balance" is untouched.
# "Logout" is untouched except it now renders a page break before its text.
$newLogoutXml = Add-PageBreakBeforeText $logoutXml "Logout"

$fullBlock = $registerXml + $newForgotPwXml + $newProfileXml + $newTransactionsXml + `
    $newDepositXml + $newWithdrawalXml + $acctBalXml + $newLogoutXml

# ---------------------------------------------------------------------------
# Replace the whole Signup..Logout run of bullet paragraphs in one shot.
# ---------------------------------------------------------------------------

$startPara = $d.Paragraphs($idxSignup)
$endPara = $d.Paragraphs($idxLogout)
$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$blockRange.InsertXML($fullBlock)

# ---------------------------------------------------------------------------
# "Nonfunctional Requirements" no longer carries the rendered page break
# (it moved up to sit in front of "Logout" instead).
# ---------------------------------------------------------------------------

$idxNonfunc2 = Find-ParaIndexByText $d "Nonfunctional Requirements"
$nonfuncPara = $d.Paragraphs($idxNonfunc2)
$newNonfuncXml = Remove-PageBreak (Get-ParaXml $nonfuncPara)
$nonfuncPara.Range.InsertXML($newNonfuncXml)

Write-Host "Edit complete"
